$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 419, pushing existing rows 419:447 down to 420:448
$ws.Rows("419:419").Insert()

# Populate the newly inserted row 419 with the new weekly record
$ws.Range("A419").Value = 4
$ws.Range("B419").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C419").Value = "Los Lagos"
$ws.Range("D419").Value = 44714
$ws.Range("E419").Value = 10
$ws.Range("F419").Value = 100112006
$ws.Range("G419").Value = "Repollo"
$ws.Range("H419").Value = "Crespo record"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 500
$ws.Range("K419").Value = 1700
$ws.Range("L419").Value = 1700
$ws.Range("M419").Value = 1700
$ws.Range("N419").Value = "$/unidad"
$ws.Range("O419").Value = "Región Metropolitana"
$ws.Range("P419").Value = 1700
$ws.Range("Q419").Value = 1
$ws.Range("R419").Value = "Hortaliza"
